$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- [CNN] Best Hyperparameters table (rows 5-9): Dropout column (O) now
# --- formatted with one decimal place ("0.0"), same for the redone CNN_it
# --- block (rows 13-17).
$ws.Range("O5").NumberFormat  = "0.0"
$ws.Range("O6").NumberFormat  = "0.0"
$ws.Range("O7").NumberFormat  = "0.0"
$ws.Range("O8").NumberFormat  = "0.0"
$ws.Range("O9").NumberFormat  = "0.0"
$ws.Range("O13").NumberFormat = "0.0"
$ws.Range("O14").NumberFormat = "0.0"
$ws.Range("O15").NumberFormat = "0.0"
$ws.Range("O16").NumberFormat = "0.0"
$ws.Range("O17").NumberFormat = "0.0"

# --- Re-plotted Y axis range tweak: MIN for the LSTM German row went from
# --- 0.65 to 0.6 so every chart shares the same Y axis range.
$ws.Range("AA6").Value = 0.6

# --- Newly obtained CNN_it results: fill in the German row (14) of the
# --- "[CNN] Best Hyperparameters, Best Epoch and Val Accuracy" /
# --- "[CNN] Test Set Results" tables.
$ws.Range("J14").Value = 2048
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 0.001
$ws.Range("O14").Value = 0.1
$ws.Range("P14").Value = 31
$ws.Range("Q14").Value = 0.87309999999999999

$ws.Range("T14").Value = 0.72860000000000003
$ws.Range("U14").Value = 0.82079999999999997
$ws.Range("V14").Value = 0.82110000000000005
$ws.Range("W14").Value = 0.81459999999999999
$ws.Range("X14").Value = 0.83679999999999999
$ws.Range("Y14").Value = 0.81200000000000006

# --- Fill in the Italian row (15): training-time summary (C/D/E) plus the
# --- hyperparameters / test-set results tables.
$ws.Range("C15").Value = "2h 06m"
$ws.Range("D15").Value = "1m 19s"
$ws.Range("E15").Value = "1.58s"

$ws.Range("J15").Value = 2048
$ws.Range("K15").Value = 7
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0.001
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 20
$ws.Range("Q15").Value = 0.88929999999999998

$ws.Range("T15").Value = 0.69189999999999996
$ws.Range("U15").Value = 0.85250000000000004
$ws.Range("V15").Value = 0.85170000000000001
$ws.Range("W15").Value = 0.86960000000000004
$ws.Range("X15").Value = 0.83420000000000005
$ws.Range("Y15").Value = 0.85119999999999996

# --- Selection moved (and the old fixed topLeftCell scroll-lock cleared)
# --- after redoing all the plots.
$ws.Range("V30").Select()
